$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 159, shifting existing rows 159-166 down to 161-168.
$ws.Rows.Item(159).Resize(2).Insert()

# Copy the style of the date column (D) from the row that is now 161 (original 159)
# into the two new rows so the new D cells keep the date format.
$ws.Range("D161").Copy()
$ws.Range("D159:D160").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 159: new Flame Seedless entry
$ws.Range("A159").Value = 11
$ws.Range("B159").Value = "Vega Monumental Concepción"
$ws.Range("C159").Value = "Bíobío"
$ws.Range("D159").Value = 44931
$ws.Range("E159").Value = 8
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100109
$ws.Range("H159").Value = "Uva"
$ws.Range("I159").Value = 100109001
$ws.Range("J159").Value = "Uva"
$ws.Range("K159").Value = "Flame Seedless"
$ws.Range("L159").Value = "Primera"
$ws.Range("M159").Value = 100
$ws.Range("N159").Value = 12000
$ws.Range("O159").Value = 13000
$ws.Range("P159").Value = 12500
$ws.Range("Q159").Value = "$/bandeja 10 kilos"
$ws.Range("R159").Value = "Provincia de Limarí"
$ws.Range("S159").Value = 1250
$ws.Range("T159").Value = 10

# Row 160: new Superior Seedless entry
$ws.Range("A160").Value = 11
$ws.Range("B160").Value = "Vega Monumental Concepción"
$ws.Range("C160").Value = "Bíobío"
$ws.Range("D160").Value = 44931
$ws.Range("E160").Value = 8
$ws.Range("F160").Value = "Fruta"
$ws.Range("G160").Value = 100109
$ws.Range("H160").Value = "Uva"
$ws.Range("I160").Value = 100109001
$ws.Range("J160").Value = "Uva"
$ws.Range("K160").Value = "Superior Seedless"
$ws.Range("L160").Value = "Primera"
$ws.Range("M160").Value = 100
$ws.Range("N160").Value = 14000
$ws.Range("O160").Value = 15000
$ws.Range("P160").Value = 14500
$ws.Range("Q160").Value = "$/bandeja 10 kilos"
$ws.Range("R160").Value = "Provincia de Limarí"
$ws.Range("S160").Value = 1450
$ws.Range("T160").Value = 10
